$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "want to go" counts in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5330
$ws1.Range("F4").Value = 11071
$ws1.Range("F8").Value = 220
$ws1.Range("F9").Value = 942

# Sheet "全部类型" (All Types) - same underlying rows, update column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5330
$ws4.Range("F7").Value = 11071
$ws4.Range("F13").Value = 220
$ws4.Range("F14").Value = 942
